# Weekly update for Haba / Vega Central Mapocho de Santiago:
# a new price-report row is inserted at row 204 (pushing the existing
# rows 204:216 down to 205:217), and the new row is populated with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 204:216 down to 205:217, leaving row 204 free.
$ws.Rows("204:204").Insert()

# Fill in the new row 204 with this week's report.
$ws.Cells.Item(204, 1).Value  = 9
$ws.Cells.Item(204, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(204, 3).Value  = "Metropolitana"
$ws.Cells.Item(204, 4).Value  = 44714
$ws.Cells.Item(204, 5).Value  = 13
$ws.Cells.Item(204, 6).Value  = 100112026
$ws.Cells.Item(204, 7).Value  = "Haba"
$ws.Cells.Item(204, 8).Value  = "Sin especificar"
$ws.Cells.Item(204, 9).Value  = "Primera"
$ws.Cells.Item(204, 10).Value = 52
$ws.Cells.Item(204, 11).Value = 18000
$ws.Cells.Item(204, 12).Value = 19000
$ws.Cells.Item(204, 13).Value = 18500
$ws.Cells.Item(204, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(204, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(204, 16).Value = 740
$ws.Cells.Item(204, 17).Value = 25
$ws.Cells.Item(204, 18).Value = "Hortaliza"
